$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - copy/paste-special constant used to replicate an
# existing cell style onto a new cell without minting duplicate style
# records (NumberFormat assignment on each cell individually creates a
# brand-new style entry every time, so formats are seeded once and then
# propagated via PasteSpecial).
$xlPasteFormats = -4122

# New rows of time-investment data (rows 8-14)
$rows = @(
    @{ Row=8;  No=6;  Tipo="Lluvia de ideas";         Fecha=43503; Inicio=0.2638888888888889; Fin=0.28125;             Interrup=10; Delta=15 },
    @{ Row=9;  No=7;  Tipo="Preguntas";               Fecha=43503; Inicio=0.28125;             Fin=0.28819444444444448; Interrup=0;  Delta=10 },
    @{ Row=10; No=8;  Tipo="Codings";                 Fecha=43506; Inicio=0.375;               Fin=0.45833333333333331; Interrup=30; Delta=90 },
    @{ Row=11; No=9;  Tipo="Mapa de Empatía";         Fecha=43508; Inicio=0.16666666666666666; Fin=0.18055555555555555; Interrup=15; Delta=5 },
    @{ Row=12; No=10; Tipo="Diagrama de caso de uso"; Fecha=43510; Inicio=0.21527777777777779; Fin=0.22916666666666666; Interrup=5;  Delta=15 },
    @{ Row=13; No=11; Tipo="Descripción de caso";     Fecha=43510; Inicio=0.22916666666666666; Fin=0.23958333333333334; Interrup=5;  Delta=10 },
    @{ Row=14; No=12; Tipo="Diagrama de actividad";   Fecha=43510; Inicio=0.23958333333333334; Fin=0.28125;             Interrup=10; Delta=50 }
)

$firstDataRow = $true

foreach ($r in $rows) {
    $i = $r.Row

    # --- values ---
    $ws.Cells.Item($i, 1).Value = $r.No          # No.
    $ws.Cells.Item($i, 2).Value = $r.Tipo         # Tipo
    $ws.Cells.Item($i, 3).Value = $r.Fecha        # Fecha
    $ws.Cells.Item($i, 4).Value = $r.Inicio       # Inicio
    $ws.Cells.Item($i, 5).Value = $r.Fin          # Fin
    $ws.Cells.Item($i, 6).Value = $r.Interrup     # Tiempo Interrupción
    $ws.Cells.Item($i, 7).Value = $r.Delta        # Tiempo Delta
    $ws.Cells.Item($i, 8).Value = "1er corte"     # Fase

    # --- plain cells (A, B, F, G, H) reuse the un-decorated style already
    #     used by the existing data rows ---
    foreach ($col in 1,2,6,7,8) {
        $ws.Range("A3").Copy()
        $ws.Cells.Item($i, $col).PasteSpecial($xlPasteFormats)
    }

    if ($firstDataRow) {
        # First new row: mint the date/time formats once.
        $ws.Cells.Item($i, 3).NumberFormat = "mm-dd-yy"
        $ws.Cells.Item($i, 4).NumberFormat = "h:mm"
        $ws.Cells.Item($i, 5).NumberFormat = "h:mm"
        $firstDataRow = $false
    } else {
        # Subsequent rows: copy the formats instead of reassigning them,
        # so the style table doesn't grow a new entry per cell.
        $ws.Range("C8").Copy()
        $ws.Cells.Item($i, 3).PasteSpecial($xlPasteFormats)
        $ws.Range("D8").Copy()
        $ws.Cells.Item($i, 4).PasteSpecial($xlPasteFormats)
        $ws.Range("D8").Copy()
        $ws.Cells.Item($i, 5).PasteSpecial($xlPasteFormats)
    }
}

# Column B is now wider to fit the longer "Tipo" labels
$ws.Columns("B").ColumnWidth = 21.8

# Selection moves to the last-edited cell
$ws.Range("G14").Select()
